$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.527.89'
$ws.Range("E2").Value = '  -1.54%  '

$ws.Range("D3").Value = '2.765.08'
$ws.Range("E3").Value = '  -2.74%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '356.95'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.43%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.82'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -4.30%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.553'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -3.09%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.582'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -4.17%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.37'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -4.33%  '

$ws.Range("E11").Value = '  +4.14%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0839'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -3.72%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.56'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.85%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.55'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.70%  '

$ws.Range("D15").Value = '3.209.44'
$ws.Range("E15").Value = '  -2.44%  '

$ws.Range("D16").Value = '2.774.83'
$ws.Range("E16").Value = '  -1.16%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.921'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.83%  '

$ws.Range("D18").Value = '51.518.49'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.60'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.11%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.06'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.64%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.05'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -3.23%  '

$ws.Range("D22").Value = '0.0₃0962'
$ws.Range("E22").Value = '  -3.87%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.70'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.29%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '267.13'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.28%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.75'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -3.14%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.14'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -3.66%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.06%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.162'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +13.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.10'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.66%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.81'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.69%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '51.74'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.48%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.04'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0437'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -9.84%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0834'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.31%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.12'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -8.31%  '

$ws.Range("E37").Value = '  +0.00%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.61'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.33%  '

$ws.Range("E39").Value = '  -5.28%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.94'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -5.67%  '

$ws.Range("E41").Value = '  -3.32%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.48'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.35%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '120.15'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -5.14%  '

$ws.Range("E44").Value = '  -3.58%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.51'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -7.24%  '

$ws.Range("D46").Value = '2.077.44'
$ws.Range("E46").Value = '  -0.89%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.23'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.87%  '

$ws.Range("E48").Value = '  +0.03%  '

$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.52'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -6.57%  '

$ws.Range("B50").Value = 'SEI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.919'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -5.31%  '

$ws.Range("E51").Value = '  +3.95%  '
